# "remove column from alcohol data"
#
# Sheet1 ("measurement-4") carries a redundant trailing column: delete
# column M so the former column N's values shift left into M, and the
# sheet's used range shrinks from A1:N119 to A1:M119.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Columns("M").Delete() | Out-Null

# Park the selection on the new last column (matches the post-edit cursor
# position recorded for this sheet).
$ws.Range("M1").Select() | Out-Null
